$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for Coliflor at Terminal
# Hortofrutícola Agro Chillán. Insert it as a new row right after the
# existing row 393, pushing the remaining historical rows (old 394..425)
# down by one (they become 395..426).
$ws.Rows.Item(394).Insert()

$ws.Cells.Item(394, 1).Value = 7
$ws.Cells.Item(394, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(394, 3).Value = "Ñuble"
$ws.Cells.Item(394, 4).Value = 45013
$ws.Cells.Item(394, 5).Value = 16
$ws.Cells.Item(394, 6).Value = 100112008
$ws.Cells.Item(394, 7).Value = "Coliflor"
$ws.Cells.Item(394, 8).Value = "Sin especificar"
$ws.Cells.Item(394, 9).Value = "Primera"
$ws.Cells.Item(394, 10).Value = 100
$ws.Cells.Item(394, 11).Value = 1200
$ws.Cells.Item(394, 12).Value = 1300
$ws.Cells.Item(394, 13).Value = 1240
$ws.Cells.Item(394, 14).Value = "$/unidad"
$ws.Cells.Item(394, 15).Value = "Región del Maule"
$ws.Cells.Item(394, 16).Value = 1240
$ws.Cells.Item(394, 17).Value = 1
$ws.Cells.Item(394, 18).Value = "Hortaliza"
